$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update financial figures to restated values
$ws.Range("D2").Value = 65150
$ws.Range("E2").Value = -2350
$ws.Range("F2").Value = -2350
$ws.Range("G2").Value = -4323
$ws.Range("H2").Value = 218
$ws.Range("I2").Value = 374
$ws.Range("J2").Value = -156
$ws.Range("K2").Value = 72658
$ws.Range("L2").Value = 65802
$ws.Range("M2").Value = 6855
$ws.Range("N2").Value = 5945
$ws.Range("O2").Value = 910
$ws.Range("P2").Value = 9114
$ws.Range("Q2").Value = 371
$ws.Range("R2").Value = 1796
$ws.Range("S2").Value = -6198
$ws.Range("T2").Value = 3584
$ws.Range("U2").Value = -3213
$ws.Range("V2").Value = 53204
$ws.Range("W2").Value = -3.61
$ws.Range("X2").Value = 0.33
$ws.Range("Y2").Value = 7.04
$ws.Range("Z2").Value = 0.27
$ws.Range("AA2").Value = 959.87
$ws.Range("AB2").Value = -26.62
$ws.Range("AC2").Value = 1292
$ws.Range("AD2").Value = 46.06
$ws.Range("AE2").Value = 19812
$ws.Range("AF2").Value = 3
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 28149732

# Row 3: update financial figures to restated values
$ws.Range("D3").Value = 56451
$ws.Range("E3").Value = -2793
$ws.Range("F3").Value = -2793
$ws.Range("G3").Value = -7053
$ws.Range("H3").Value = -6805
$ws.Range("I3").Value = -6695
$ws.Range("J3").Value = -109
$ws.Range("K3").Value = 59425
$ws.Range("L3").Value = 57139
$ws.Range("M3").Value = 2286
$ws.Range("N3").Value = 1857
$ws.Range("O3").Value = 430
$ws.Range("P3").Value = 11825
$ws.Range("Q3").Value = -1566
$ws.Range("R3").Value = 1770
$ws.Range("S3").Value = -2637
$ws.Range("T3").Value = 1406
$ws.Range("U3").Value = -2972
$ws.Range("V3").Value = 49036
$ws.Range("W3").Value = -4.95
$ws.Range("X3").Value = -12.05
$ws.Range("Y3").Value = -171.63
$ws.Range("Z3").Value = -10.3
$ws.Range("AA3").Value = 2499.12
$ws.Range("AB3").Value = -73.23999999999999
$ws.Range("AC3").Value = -19993
$ws.Range("AD3").Value = -1.29
$ws.Range("AE3").Value = 5039
$ws.Range("AF3").Value = 5.12
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 34989656

# Row 4: update financial figures to restated values
$ws.Range("D4").Value = 45848
$ws.Range("E4").Value = -8334
$ws.Range("F4").Value = -8334
$ws.Range("G4").Value = -4155
$ws.Range("H4").Value = -4582
$ws.Range("I4").Value = -4596
$ws.Range("J4").Value = 14
$ws.Range("K4").Value = 43981
$ws.Range("L4").Value = 34468
$ws.Range("M4").Value = 9513
$ws.Range("N4").Value = 9478
$ws.Range("O4").Value = 35
$ws.Range("P4").Value = 8986
$ws.Range("Q4").Value = -6390
$ws.Range("R4").Value = 14739
$ws.Range("S4").Value = -5030
$ws.Range("T4").Value = 1552
$ws.Range("U4").Value = -7942
$ws.Range("V4").Value = 26926
$ws.Range("W4").Value = -18.18
$ws.Range("X4").Value = -9.99
$ws.Range("Y4").Value = -81.09999999999999
$ws.Range("Z4").Value = -8.859999999999999
$ws.Range("AA4").Value = 362.31
$ws.Range("AB4").Value = 12.03
$ws.Range("AC4").Value = -4123
$ws.Range("AD4").Value = -1.46
$ws.Range("AE4").Value = 4656
$ws.Range("AF4").Value = 1.3
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 201734167

# Row 5: update financial figures to restated values
$ws.Range("D5").Value = 50280
$ws.Range("E5").Value = -4068
$ws.Range("F5").Value = -4068
$ws.Range("G5").Value = -11870
$ws.Range("H5").Value = -11907
$ws.Range("I5").Value = -11910
$ws.Range("J5").Value = 3
$ws.Range("K5").Value = 36024
$ws.Range("L5").Value = 27055
$ws.Range("M5").Value = 8969
$ws.Range("N5").Value = 8948
$ws.Range("O5").Value = 21
$ws.Range("P5").Value = 15683
$ws.Range("Q5").Value = -2754
$ws.Range("R5").Value = -2280
$ws.Range("S5").Value = 6615
$ws.Range("T5").Value = 3025
$ws.Range("U5").Value = -5779
$ws.Range("V5").Value = 20520
$ws.Range("W5").Value = -8.09
$ws.Range("X5").Value = -23.68
$ws.Range("Y5").Value = -129.27
$ws.Range("Z5").Value = -29.77
$ws.Range("AA5").Value = 301.64
$ws.Range("AB5").Value = -70.55
$ws.Range("AC5").Value = -5475
$ws.Range("AD5").Value = -0.92
$ws.Range("AE5").Value = 2853
$ws.Range("AF5").Value = 1.76
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 311805827

# Row 6: update financial figures to restated values
$ws.Range("D6").Value = 52221
$ws.Range("E6").Value = -5587
$ws.Range("F6").Value = -5587
$ws.Range("G6").Value = -7789
$ws.Range("H6").Value = -7906
$ws.Range("I6").Value = -7907
$ws.Range("K6").Value = 41214
$ws.Range("L6").Value = 30818
$ws.Range("M6").Value = 10397
$ws.Range("N6").Value = 10385
$ws.Range("P6").Value = 15783
$ws.Range("Q6").Value = -4169
$ws.Range("R6").Value = -6728
$ws.Range("S6").Value = 9874
$ws.Range("T6").Value = 3201
$ws.Range("U6").Value = -7370
$ws.Range("V6").Value = 23461
$ws.Range("W6").Value = -10.7
$ws.Range("X6").Value = -15.14
$ws.Range("Y6").Value = -81.8
$ws.Range("Z6").Value = -20.47
$ws.Range("AA6").Value = 296.42
$ws.Range("AB6").Value = -125.3
$ws.Range("AC6").Value = -2521
$ws.Range("AD6").Value = -1.47
$ws.Range("AE6").Value = 3290
$ws.Range("AF6").Value = 1.12
$ws.Range("AG6").Value = 0
$ws.Range("AJ6").Value = 313805827

# Row 6: AH/AI no longer reported for this period
$ws.Range("AH6").ClearContents()
$ws.Range("AI6").ClearContents()

# Row 7: estimates removed, only identifying columns remain
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AC7").ClearContents()
$ws.Range("AD7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

# Row 8: estimates removed, only identifying columns remain
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AC8").ClearContents()
$ws.Range("AD8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9: estimates removed, only identifying columns remain
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()
